# Generate Report for Handoff
#
# Replaces the previous "handed back" localization-status snapshot with a
# fresh "ready for handoff" snapshot: new source-file UUIDs, new status
# text, new handoff timestamps, and clears out the now-stale
# "Latest Target File" / "Latest Handback File" columns (F/G) on the
# per-language sheets (their hyperlinks are removed too).

$wb = $excel.ActiveWorkbook

# ---- New data values ---------------------------------------------------
$newMdA      = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.md"
$newMdB      = "fffff5089c33-728a-49b8-bee3-dc59b445d663.md"
$status      = "Ready for handoff"
$overviewDt  = "2016-29-18 07:29:06"

$zhXlf       = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.zh-cn.xlf"
$zhHandoffDt = "2016-03-18 07:29:03"

$deXlf       = "6ed8456b-c2ad-4330-9572-caf6d7fa80c0.41d923ced29b8f299034b77d0df713481a4d1485.de-de.xlf"
$deHandoffDt = "2016-03-18 07:29:06"

$emptyHandback = "0001-01-01 00:00:00"

# Helper: update the display text of an existing hyperlink in place
# (keeps its address / r:id / cell style untouched).
function Set-HyperlinkDisplay {
    param($ws, [string]$addr, [string]$newDisplay)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $newDisplay
            return
        }
    }
}

# Helper: remove a single hyperlink (identified by its cell address)
# without disturbing the others.
function Remove-HyperlinkAt {
    param($ws, [string]$addr)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            return
        }
    }
}

# =========================================================================
# Sheet "Overview"
# =========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newMdA
$ws.Range("B2").Value = $status
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $overviewDt

$ws.Range("A3").Value = $newMdB
$ws.Range("B3").Value = $status
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $overviewDt

Set-HyperlinkDisplay $ws '$A$2' $newMdA
Set-HyperlinkDisplay $ws '$A$3' $newMdB

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newMdA
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $zhXlf
$ws.Range("E2").Value = $zhHandoffDt
$ws.Range("H2").Value = $emptyHandback
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $newMdB
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $zhXlf
$ws.Range("E3").Value = $zhHandoffDt
$ws.Range("H3").Value = $emptyHandback
$ws.Range("I3").Value = "Include"

Remove-HyperlinkAt $ws '$F$2'
Remove-HyperlinkAt $ws '$G$2'
Remove-HyperlinkAt $ws '$F$3'
Remove-HyperlinkAt $ws '$G$3'

$ws.Range("F2:G3").Clear()

Set-HyperlinkDisplay $ws '$A$2' $newMdA
Set-HyperlinkDisplay $ws '$B$2' ".md"
Set-HyperlinkDisplay $ws '$D$2' $zhXlf
Set-HyperlinkDisplay $ws '$A$3' $newMdB
Set-HyperlinkDisplay $ws '$B$3' ".md"
Set-HyperlinkDisplay $ws '$D$3' $zhXlf

# =========================================================================
# Sheet "de-de"
# =========================================================================
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newMdA
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $status
$ws.Range("D2").Value = $deXlf
$ws.Range("E2").Value = $deHandoffDt
$ws.Range("H2").Value = $emptyHandback
$ws.Range("I2").Value = "Include"

$ws.Range("A3").Value = $newMdB
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $status
$ws.Range("D3").Value = $deXlf
$ws.Range("E3").Value = $deHandoffDt
$ws.Range("H3").Value = $emptyHandback
$ws.Range("I3").Value = "Include"

Remove-HyperlinkAt $ws '$F$2'
Remove-HyperlinkAt $ws '$G$2'
Remove-HyperlinkAt $ws '$F$3'
Remove-HyperlinkAt $ws '$G$3'

$ws.Range("F2:G3").Clear()

Set-HyperlinkDisplay $ws '$A$2' $newMdA
Set-HyperlinkDisplay $ws '$B$2' ".md"
Set-HyperlinkDisplay $ws '$D$2' $deXlf
Set-HyperlinkDisplay $ws '$A$3' $newMdB
Set-HyperlinkDisplay $ws '$B$3' ".md"
Set-HyperlinkDisplay $ws '$D$3' $deXlf

Write-Host "Report regenerated for handoff"
